$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for season record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style from an existing header cell (e.g. AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Fill season record for every data row (2-49) with the team's 2000 record
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 94   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 68   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
